$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder word-list with the "Testcodeword" codeword list
for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i, 1).Value = "Testcodeword$i"
}

# Autofit column A to the new (wider) content
$ws.Columns.Item(1).AutoFit() | Out-Null

# Reflect the new selection (A11:A14, active cell A14) used to show the ack status
$ws.Range("A11:A14").Select() | Out-Null
$ws.Range("A14").Activate() | Out-Null
